$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the mailto: hyperlink that lives on the Email cell (B2) before the
# column is removed.
$ws.Range("B2").Hyperlinks.Delete()

# Remove the "Email" column (column B) entirely - this removes the
# hyperlink cell/style along with its data and shifts Phone/License
# Number/License Image one column to the left.
$ws.Columns("B").Delete()

# Restore the selection to match the post-edit workbook (column B selected).
$ws.Range("B1:B1048576").Select()

$wb.Save()
